$d = $word.ActiveDocument

$pairs = @(
    @("638÷3=212, 2", "595÷9=66, 1"),
    @("227÷5=45, 2", "679÷5=135, 4"),
    @("659÷5=131, 4", "336÷2=168, 0"),
    @("387÷6=64, 3", "997÷2=498, 1"),
    @("957÷6=159, 3", "618÷4=154, 2"),
    @("256÷2=128, 0", "662÷5=132, 2"),
    @("746÷9=82, 8", "331÷4=82, 3"),
    @("801÷2=400, 1", "586÷6=97, 4"),
    @("120÷2=60, 0", "439÷6=73, 1"),
    @("586÷7=83, 5", "869÷7=124, 1"),
    @("898÷3=299, 1", "457÷2=228, 1"),
    @("730÷7=104, 2", "428÷8=53, 4"),
    @("801÷7=114, 3", "143÷8=17, 7"),
    @("670÷6=111, 4", "409÷5=81, 4"),
    @("724÷7=103, 3", "595÷6=99, 1"),
    @("948÷9=105, 3", "794÷3=264, 2"),
    @("502÷2=251, 0", "671÷8=83, 7"),
    @("184÷6=30, 4", "591÷3=197, 0"),
    @("676÷6=112, 4", "210÷6=35, 0"),
    @("390÷3=130, 0", "371÷5=74, 1"),
    @("733÷7=104, 5", "957÷9=106, 3"),
    @("771÷5=154, 1", "244÷2=122, 0"),
    @("217÷4=54, 1", "261÷2=130, 1"),
    @("733÷6=122, 1", "126÷8=15, 6"),
    @("108÷6=18, 0", "666÷5=133, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
